# Scheduled-runner refresh: push newly-fetched market-board averages (and
# recomputed leve profit figures) into each job sheet's data rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1667283.4
$ws.Range("I2").Value = 5000100.5
$ws.Range("K2").Value = 5000100.5
$ws.Range("M2").Value = -4999987.5

$ws.Range("H15").Value = 10527908
$ws.Range("I15").Value = 10527908
$ws.Range("K15").Value = 31583724
$ws.Range("M15").Value = -31583555

$ws.Range("H113").Value = 2622.9583
$ws.Range("I113").Value = 2943.2666
$ws.Range("J113").Value = 2089.111
$ws.Range("K113").Value = 2943.2666
$ws.Range("L113").Value = 2089.111
$ws.Range("M113").Value = 310.7334000000001
$ws.Range("N113").Value = -8597.111000000001

$ws.Range("H116").Value = 4454.1
$ws.Range("I116").Value = 4872.5
$ws.Range("K116").Value = 4872.5
$ws.Range("M116").Value = -1430.5

$ws.Range("H132").Value = 223751.69
$ws.Range("I132").Value = 1475.9445
$ws.Range("J132").Value = 1112854.6
$ws.Range("K132").Value = 4427.833500000001
$ws.Range("L132").Value = 3338563.8
$ws.Range("M132").Value = -1897.833500000001
$ws.Range("N132").Value = -3343623.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2099.6667
$ws.Range("I45").Value = 1737.125
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 1737.125
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -1360.125
$ws.Range("N45").Value = -5754

$ws.Range("H61").Value = 1405.3948
$ws.Range("I61").Value = 1242.3214
$ws.Range("J61").Value = 1862
$ws.Range("K61").Value = 1242.3214
$ws.Range("L61").Value = 1862
$ws.Range("M61").Value = -1030.3214
$ws.Range("N61").Value = -2286

$ws.Range("H63").Value = 1586.25
$ws.Range("I63").Value = 1636
$ws.Range("J63").Value = 1238
$ws.Range("K63").Value = 1636
$ws.Range("L63").Value = 1238
$ws.Range("M63").Value = -950
$ws.Range("N63").Value = -2610

$ws.Range("H66").Value = 1586.25
$ws.Range("I66").Value = 1636
$ws.Range("J66").Value = 1238
$ws.Range("K66").Value = 8180
$ws.Range("L66").Value = 6190
$ws.Range("M66").Value = -4748
$ws.Range("N66").Value = -13054

$ws.Range("H74").Value = 818.7857
$ws.Range("I74").Value = 769.8043
$ws.Range("J74").Value = 1044.1
$ws.Range("K74").Value = 769.8043
$ws.Range("L74").Value = 1044.1
$ws.Range("M74").Value = 104.1957
$ws.Range("N74").Value = -2792.1

$ws.Range("H77").Value = 818.7857
$ws.Range("I77").Value = 769.8043
$ws.Range("J77").Value = 1044.1
$ws.Range("K77").Value = 3849.0215
$ws.Range("L77").Value = 5220.5
$ws.Range("M77").Value = 518.9785000000002
$ws.Range("N77").Value = -13956.5

$ws.Range("H132").Value = 1610.871
$ws.Range("I132").Value = 1226.7693
$ws.Range("J132").Value = 1888.2778
$ws.Range("K132").Value = 3680.3079
$ws.Range("L132").Value = 5664.8334
$ws.Range("M132").Value = -1150.3079
$ws.Range("N132").Value = -10724.8334

$ws.Range("H136").Value = 1405.3948
$ws.Range("I136").Value = 1242.3214
$ws.Range("J136").Value = 1862
$ws.Range("K136").Value = 3726.9642
$ws.Range("L136").Value = 5586
$ws.Range("M136").Value = -1176.9642
$ws.Range("N136").Value = -10686

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1368.6666
$ws.Range("I94").Value = 919.6667
$ws.Range("J94").Value = 2266.6667
$ws.Range("K94").Value = 919.6667
$ws.Range("L94").Value = 2266.6667
$ws.Range("M94").Value = -468.6667
$ws.Range("N94").Value = -3168.6667

$ws.Range("H105").Value = 4116.7915
$ws.Range("I105").Value = 3657.2856
$ws.Range("J105").Value = 7333.3335
$ws.Range("K105").Value = 3657.2856
$ws.Range("L105").Value = 7333.3335
$ws.Range("M105").Value = -1910.2856
$ws.Range("N105").Value = -10827.3335

$ws.Range("H107").Value = 11762.75
$ws.Range("I107").Value = 765.3
$ws.Range("J107").Value = 66750
$ws.Range("K107").Value = 765.3
$ws.Range("L107").Value = 66750
$ws.Range("M107").Value = 1154.7
$ws.Range("N107").Value = -70590

$ws.Range("H133").Value = 66526.664
$ws.Range("J133").Value = 66526.664
$ws.Range("L133").Value = 66526.664
$ws.Range("N133").Value = -76646.664

$ws.Range("H134").Value = 48656.066
$ws.Range("I134").Value = 3494.5
$ws.Range("J134").Value = 102849.95
$ws.Range("K134").Value = 10483.5
$ws.Range("L134").Value = 308549.85
$ws.Range("M134").Value = -7948.5
$ws.Range("N134").Value = -313619.85

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 295.0303
$ws.Range("I107").Value = 143.53847
$ws.Range("J107").Value = 393.5
$ws.Range("K107").Value = 143.53847
$ws.Range("L107").Value = 393.5
$ws.Range("M107").Value = 1776.46153
$ws.Range("N107").Value = -4233.5

$ws.Range("H122").Value = 1667386.9
$ws.Range("I122").Value = 2000764.6
$ws.Range("J122").Value = 498
$ws.Range("K122").Value = 6002293.800000001
$ws.Range("L122").Value = 1494
$ws.Range("M122").Value = -5999843.800000001
$ws.Range("N122").Value = -6394

$ws.Range("H132").Value = 2445.25
$ws.Range("I132").Value = 1786.6522
$ws.Range("J132").Value = 4128.3335
$ws.Range("K132").Value = 5359.9566
$ws.Range("L132").Value = 12385.0005
$ws.Range("M132").Value = -2829.9566
$ws.Range("N132").Value = -17445.0005

$ws.Range("H134").Value = 3262
$ws.Range("I134").Value = 2938.8462
$ws.Range("J134").Value = 3787.125
$ws.Range("K134").Value = 8816.5386
$ws.Range("L134").Value = 11361.375
$ws.Range("M134").Value = -6281.5386
$ws.Range("N134").Value = -16431.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2718.3333
$ws.Range("I136").Value = 1968.8889
$ws.Range("J136").Value = 4966.6665
$ws.Range("K136").Value = 5906.6667
$ws.Range("L136").Value = 14899.9995
$ws.Range("M136").Value = -806.6666999999998
$ws.Range("N136").Value = -25099.9995

$ws.Range("H140").Value = 5044.6665
$ws.Range("I140").Value = 5963.3335
$ws.Range("J140").Value = 3207.3333
$ws.Range("K140").Value = 17890.0005
$ws.Range("L140").Value = 9621.999899999999
$ws.Range("M140").Value = -12710.0005
$ws.Range("N140").Value = -19981.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4376.4062
$ws.Range("I70").Value = 4083.96
$ws.Range("K70").Value = 4083.96
$ws.Range("M70").Value = -3813.96

$ws.Range("H73").Value = 4376.4062
$ws.Range("I73").Value = 4083.96
$ws.Range("K73").Value = 4083.96
$ws.Range("M73").Value = -3147.96

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 461.2857
$ws.Range("I55").Value = 493.33334
$ws.Range("K55").Value = 493.33334
$ws.Range("M55").Value = -320.33334

$ws.Range("H61").Value = 6846.409
$ws.Range("I61").Value = 7137.9473
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 7137.9473
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -6935.9473
$ws.Range("N61").Value = -5404

$ws.Range("H113").Value = 6846.409
$ws.Range("I113").Value = 7137.9473
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 7137.9473
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -4967.9473
$ws.Range("N113").Value = -9340

$ws.Range("H136").Value = 3122.524
$ws.Range("I136").Value = 1345.9474
$ws.Range("J136").Value = 20000
$ws.Range("K136").Value = 4037.8422
$ws.Range("L136").Value = 60000
$ws.Range("M136").Value = -1487.8422
$ws.Range("N136").Value = -65100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 408.53333
$ws.Range("I126").Value = 398.92307
$ws.Range("J126").Value = 471
$ws.Range("K126").Value = 1196.76921
$ws.Range("L126").Value = 1413
$ws.Range("M126").Value = 1273.23079
$ws.Range("N126").Value = -6353

$ws.Range("H132").Value = 1775.9474
$ws.Range("I132").Value = 1436.2667
$ws.Range("J132").Value = 3049.75
$ws.Range("K132").Value = 4308.800099999999
$ws.Range("L132").Value = 9149.25
$ws.Range("M132").Value = -1778.800099999999
$ws.Range("N132").Value = -14209.25

$ws.Range("H136").Value = 1313.317
$ws.Range("I136").Value = 1278.1025
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 3834.3075
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -1284.3075
$ws.Range("N136").Value = -11100
